# Update "想去人数" (want-to-go count) values in column F across the
# 展览 (sheet1), 演出 (sheet2), 本地生活 (sheet3) and 全部类型 (sheet4)
# worksheets, matching the data refresh captured in the diff.

$wb = $excel.ActiveWorkbook

$sheetExhibitions = $wb.Worksheets.Item("展览")
$sheetPerformances = $wb.Worksheets.Item("演出")
$sheetLocalLife = $wb.Worksheets.Item("本地生活")
$sheetAllTypes = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$sheetExhibitions.Range("F2").Value = 568
$sheetExhibitions.Range("F4").Value = 23
$sheetExhibitions.Range("F5").Value = 724
$sheetExhibitions.Range("F6").Value = 346
$sheetExhibitions.Range("F8").Value = 136
$sheetExhibitions.Range("F10").Value = 203
$sheetExhibitions.Range("F11").Value = 5772
$sheetExhibitions.Range("F12").Value = 41
$sheetExhibitions.Range("F16").Value = 540
$sheetExhibitions.Range("F17").Value = 335
$sheetExhibitions.Range("F22").Value = 88
$sheetExhibitions.Range("F26").Value = 1746
$sheetExhibitions.Range("F27").Value = 447

# 演出 (sheet2)
$sheetPerformances.Range("F2").Value = 657
$sheetPerformances.Range("F5").Value = 260
$sheetPerformances.Range("F6").Value = 293
$sheetPerformances.Range("F10").Value = 135

# 本地生活 (sheet3)
$sheetLocalLife.Range("F2").Value = 191

# 全部类型 (sheet4)
$sheetAllTypes.Range("F2").Value = 191
$sheetAllTypes.Range("F3").Value = 568
$sheetAllTypes.Range("F5").Value = 23
$sheetAllTypes.Range("F6").Value = 724
$sheetAllTypes.Range("F7").Value = 657
$sheetAllTypes.Range("F8").Value = 346
$sheetAllTypes.Range("F10").Value = 136
$sheetAllTypes.Range("F12").Value = 203
$sheetAllTypes.Range("F13").Value = 5772
$sheetAllTypes.Range("F14").Value = 41
$sheetAllTypes.Range("F19").Value = 540
$sheetAllTypes.Range("F20").Value = 335
$sheetAllTypes.Range("F25").Value = 260
$sheetAllTypes.Range("F26").Value = 293
$sheetAllTypes.Range("F31").Value = 135
$sheetAllTypes.Range("F32").Value = 88
$sheetAllTypes.Range("F36").Value = 1746
$sheetAllTypes.Range("F37").Value = 447
